# Apply updated crypto price/volume data to worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.643.08"
$ws.Range("E2").Value = "'  +3.49%  "
$ws.Range("D3").Value = "'1.859.11"
$ws.Range("E3").Value = "'  +2.20%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("D5").Value = "'273.21"
$ws.Range("E5").Value = "'  -2.05%  "
$ws.Range("D6").Value = "'0.9997"
$ws.Range("E6").Value = "'  -0.12%  "
$ws.Range("D7").Value = "'0.5272"
$ws.Range("E7").Value = "'  +3.46%  "
$ws.Range("D8").Value = "'0.3384"
$ws.Range("E8").Value = "'  -4.40%  "
$ws.Range("D9").Value = "'0.06787"
$ws.Range("E9").Value = "'  +1.44%  "
$ws.Range("E10").Value = "'  -0.70%  "
$ws.Range("D11").Value = "'0.7919"
$ws.Range("E11").Value = "'  -4.38%  "
$ws.Range("D12").Value = "'0.07741"
$ws.Range("E12").Value = "'  -1.63%  "
$ws.Range("D13").Value = "'1.871.33"
$ws.Range("E13").Value = "'  +2.85%  "
$ws.Range("D14").Value = "'89.59"
$ws.Range("E14").Value = "'  +1.99%  "
$ws.Range("D15").Value = "'5.118"
$ws.Range("E15").Value = "'  +0.64%  "
$ws.Range("D16").Value = "'0.9996"
$ws.Range("E16").Value = "'  -0.12%  "
$ws.Range("D17").Value = "'14.41"
$ws.Range("E17").Value = "'  +2.13%  "
$ws.Range("D18").Value = "'0.000007996"
$ws.Range("E18").Value = "'  -0.63%  "
$ws.Range("D19").Value = "'0.9999"
$ws.Range("D20").Value = "'26.669.07"
$ws.Range("E20").Value = "'  +3.41%  "
$ws.Range("D21").Value = "'2.109.11"
$ws.Range("E21").Value = "'  +3.08%  "
$ws.Range("D22").Value = "'4.709"
$ws.Range("E22").Value = "'  -1.07%  "
$ws.Range("D23").Value = "'9.948"
$ws.Range("E23").Value = "'  -0.59%  "
$ws.Range("D24").Value = "'6.084"
$ws.Range("E24").Value = "'  -0.49%  "
$ws.Range("D25").Value = "'2.359"
$ws.Range("E25").Value = "'  +4.89%  "
$ws.Range("D26").Value = "'145.85"
$ws.Range("E26").Value = "'  +2.43%  "
$ws.Range("D27").Value = "'1.651"
$ws.Range("E27").Value = "'  -1.21%  "
$ws.Range("D28").Value = "'17.18"
$ws.Range("E28").Value = "'  +0.14%  "
$ws.Range("D29").Value = "'112.12"
$ws.Range("E29").Value = "'  +2.46%  "
$ws.Range("D30").Value = "'4.322"
$ws.Range("E30").Value = "'  -0.52%  "
$ws.Range("D31").Value = "'4.302"
$ws.Range("E31").Value = "'  +1.47%  "
$ws.Range("D32").Value = "'0.08867"
$ws.Range("E32").Value = "'  +1.07%  "
$ws.Range("D33").Value = "'0.04907"
$ws.Range("E33").Value = "'  +0.17%  "
$ws.Range("E34").Value = "'  +1.68%  "
$ws.Range("D35").Value = "'0.7257"
$ws.Range("E35").Value = "'  -0.66%  "
$ws.Range("D36").Value = "'2.877"
$ws.Range("E36").Value = "'  -0.36%  "
$ws.Range("D37").Value = "'3.226"
$ws.Range("E37").Value = "'  +2.06%  "
$ws.Range("D38").Value = "'2.321"
$ws.Range("E38").Value = "'  -1.76%  "
$ws.Range("D39").Value = "'0.01845"
$ws.Range("E39").Value = "'  -0.64%  "
$ws.Range("D40").Value = "'0.5085"
$ws.Range("E40").Value = "'  -1.76%  "
$ws.Range("E41").Value = "'  -2.90%  "
$ws.Range("D42").Value = "'116.06"
$ws.Range("E42").Value = "'  +1.35%  "
$ws.Range("D43").Value = "'6.123"
$ws.Range("E43").Value = "'  -1.67%  "
$ws.Range("D44").Value = "'7.988"
$ws.Range("E44").Value = "'  -0.44%  "
$ws.Range("D45").Value = "'0.9990"
$ws.Range("E45").Value = "'  -0.17%  "
$ws.Range("D46").Value = "'0.4402"
$ws.Range("E46").Value = "'  -3.14%  "
$ws.Range("E47").Value = "'  -3.41%  "
$ws.Range("D48").Value = "'9.265"
$ws.Range("E48").Value = "'  +0.74%  "
$ws.Range("D49").Value = "'36.06"
$ws.Range("E49").Value = "'  -1.28%  "
$ws.Range("D50").Value = "'0.05931"
$ws.Range("E50").Value = "'  +1.53%  "
$ws.Range("D51").Value = "'1.471"
$ws.Range("E51").Value = "'  -2.02%  "
